$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (344-357), covering 2021-08-10 through 2021-08-23
$data = @(
    @(44418, 2, 14, 571.1954304365565),
    @(44419, 0, 14, 571.1954304365565),
    @(44420, 0, 13, 530.3957568339454),
    @(44421, 2, 14, 571.1954304365565),
    @(44422, 3, 11, 448.796409628723),
    @(44423, 2, 10, 407.9967360261118),
    @(44424, 1, 10, 407.9967360261118),
    @(44425, 0, 8, 326.3973888208894),
    @(44426, 0, 8, 326.3973888208894),
    @(44427, 2, 10, 407.9967360261118),
    @(44428, 1, 9, 367.1970624235006),
    @(44429, 0, 6, 244.798041615667),
    @(44430, 0, 4, 163.1986944104447),
    @(44431, 2, 5, 203.9983680130559)
)

$startRow = 344
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$endRow = $startRow + $data.Count - 1

# Copy the date-column style (column A) from the last pre-existing row so the
# new date cells match the formatting (centered, bordered, date number format)
$ws.Range("A343").Copy()
$ws.Range("A344:A$endRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false
